# "merged input from Benoit"
#
# Slide 4 ("Draft Status"): the first bullet is split in two -
#   - the old "Received and addressed comments ..." is shortened to
#     "Received comments ..."
#   - a brand-new bold-lead-in bullet "Addressed all open issues ...and
#     double-checked the IANA consideration section with the IPFIX doctors."
#     (which used to live on slide 5) is inserted right after it.
#
# Slide 5 ("Next Steps"): the now-duplicated "The authors addressed all open
# issues ..." bullet is removed, and the placeholder's autofit line-spacing
# reduction is cleared (PowerPoint recalculates it now that a bullet is gone).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4 - Content Placeholder 2
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$body4 = $slide4.Shapes.Item(2).TextFrame.TextRange

# Replace paragraph 1 in one shot (delete + insert-before) so PowerPoint
# writes a single clean run instead of diffing against the old text.
[void]$body4.Paragraphs(1).Delete()
[void]$body4.InsertBefore("Received comments from SPRING, OPSAWG and other network operators.`rAddressed all open issues and double-checked the IANA consideration section with the IPFIX doctors.`r")

# Bold just the "Addressed all open issues " lead-in of the new 2nd paragraph.
$newPara = $body4.Paragraphs(2)
$leadIn = $body4.Characters($newPara.Start, 26)
$leadIn.Font.Bold = $true

# ---------------------------------------------------------------------
# Slide 5 - Content Placeholder 2
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$shape5 = $slide5.Shapes.Item(2)
$body5 = $shape5.TextFrame.TextRange

# Drop the bullet that was merged into slide 4 above.
[void]$body5.Paragraphs(4).Delete()

# Clear the autofit line-space reduction (normAutofit lnSpcReduction="10000"
# -> normAutofit) now that the text fits without shrinking.
$shape5.TextFrame.AutoSize = 2
